# Refresh the "想去人数" (interested-count) figures in column F across the
# three populated sheets, matching the upstream data-generation run at
# commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value  = 1203
$ws.Range("F6").Value  = 79
$ws.Range("F7").Value  = 4482
$ws.Range("F8").Value  = 2656
$ws.Range("F9").Value  = 63
$ws.Range("F10").Value = 2610
$ws.Range("F15").Value = 692
$ws.Range("F16").Value = 142
$ws.Range("F17").Value = 161
$ws.Range("F19").Value = 32
$ws.Range("F22").Value = 43
$ws.Range("F24").Value = 33
$ws.Range("F26").Value = 584
$ws.Range("F27").Value = 713
$ws.Range("F28").Value = 123
$ws.Range("F29").Value = 13
$ws.Range("F30").Value = 448
$ws.Range("F32").Value = 1210
$ws.Range("F33").Value = 214
$ws.Range("F34").Value = 27
$ws.Range("F35").Value = 1283
$ws.Range("F36").Value = 2148
$ws.Range("F37").Value = 321
$ws.Range("F39").Value = 563
$ws.Range("F41").Value = 37
$ws.Range("F43").Value = 701
$ws.Range("F44").Value = 1381
$ws.Range("F45").Value = 138
$ws.Range("F47").Value = 455
$ws.Range("F48").Value = 62
$ws.Range("F49").Value = 85

# --- Sheet "演出" ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 13
$ws.Range("F4").Value = 2
$ws.Range("F5").Value = 76

# --- Sheet "全部类型" ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value  = 13
$ws.Range("F4").Value  = 79
$ws.Range("F5").Value  = 4482
$ws.Range("F6").Value  = 2656
$ws.Range("F7").Value  = 2610
$ws.Range("F10").Value = 2
$ws.Range("F11").Value = 692
$ws.Range("F12").Value = 142
$ws.Range("F13").Value = 161
$ws.Range("F15").Value = 32
$ws.Range("F18").Value = 43
$ws.Range("F20").Value = 33
$ws.Range("F21").Value = 584
$ws.Range("F22").Value = 713
$ws.Range("F23").Value = 123
$ws.Range("F24").Value = 76
$ws.Range("F27").Value = 13
$ws.Range("F28").Value = 448
$ws.Range("F30").Value = 1210
$ws.Range("F31").Value = 214
$ws.Range("F34").Value = 2148
$ws.Range("F35").Value = 321
$ws.Range("F39").Value = 563
$ws.Range("F41").Value = 37
$ws.Range("F43").Value = 701
$ws.Range("F44").Value = 1381
$ws.Range("F46").Value = 138
$ws.Range("F47").Value = 455
$ws.Range("F48").Value = 85
